$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-12-10 18:35:49"

for ($r = 2; $r -le 16; $r++) {
    $ws.Cells.Item($r, 1).Value = $newTimestamp
}
